# Test-Data.xlsx: rename the Check-in/Check-out headers, drop the
# now-unused trailing blank column D, and reset the sheet selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text fixes: "Check-in" -> "CheckIn", "Check-out" -> "CheckOut"
$ws.Range("B1").Value = "CheckIn"
$ws.Range("C1").Value = "CheckOut"

# Column D was only ever blank placeholder cells - remove it entirely.
$ws.Columns("D").Delete()

# Reset the active selection back to A1.
$ws.Range("A1").Select() | Out-Null
